$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 11 - RedesNeuronales3
$ws.Range("B11").Value = 0.8242105263157895
$ws.Range("C11").Value = 0.7674876618243893
$ws.Range("D11").Value = 0.3308457711442786
$ws.Range("E11").Value = 0.1341789052069426
$ws.Range("F11").Value = 0.7591655481289358
$ws.Range("G11").Value = 0.796
$ws.Range("H11").Value = 0.3486238532110092

# Row 12 - RedesNeuronales2
$ws.Range("B12").Value = 0.8215789473684211
$ws.Range("C12").Value = 0.7621787590751185
$ws.Range("D12").Value = 0.3407960199004975
$ws.Range("E12").Value = 0.1348464619492657
$ws.Range("F12").Value = 0.7534638409386557
$ws.Range("G12").Value = 0.796
$ws.Range("H12").Value = 0.3577981651376147

# Row 13 - RedesNeuronales
$ws.Range("B13").Value = 0.8278947368421052
$ws.Range("C13").Value = 0.7006639034467184
$ws.Range("D13").Value = 0.5199004975124378
$ws.Range("E13").Value = 0.0787716955941255
$ws.Range("F13").Value = 0.6819971899339645
$ws.Range("G13").Value = 0.792
$ws.Range("H13").Value = 0.5596330275229358
